$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B2:K11 with latest model run results (row 1 header and column A dates are unchanged) ---
$ws.Range("B2").Value = 0.35746311125172947
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("B3").Value = 0.3275176712894272
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.0038768269880811042
$ws.Range("E3").Value = 0.00003274908276443072
$ws.Range("F3").Value = -0.0007371272782482466
$ws.Range("G3").Value = 0.0006081699062767512
$ws.Range("H3").Value = -0.00004305096468067913
$ws.Range("I3").Value = -0.0012041227360983437
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.002390345643339231

$ws.Range("B4").Value = 0.30633046272307024
$ws.Range("C4").Value = -0.001366447076493017
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.0000653818879796283
$ws.Range("F4").Value = 0.000020639244225066247
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.00010437213371903527
$ws.Range("I4").Value = -0.0017103684435274048
$ws.Range("J4").Value = -0.00002124466794271749
$ws.Range("K4").Value = 0.00032420098526608765

$ws.Range("B5").Value = 0.3016693565432873
$ws.Range("C5").Value = 0.005504259730728322
$ws.Range("D5").Value = -0.0061073149052911025
$ws.Range("E5").Value = -0.00000653951679275918
$ws.Range("F5").Value = -0.0005995352634100171
$ws.Range("G5").Value = -0.0014981345521955592
$ws.Range("H5").Value = -0.00006000096147438754
$ws.Range("I5").Value = -0.0003978587733951075
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.0009569762578217622

$ws.Range("B6").Value = 0.34911972882675274
$ws.Range("C6").Value = 0.027155041000261303
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.00017364690473465342
$ws.Range("F6").Value = -0.00003874456959857522
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.00007671807285656278
$ws.Range("I6").Value = -0.0028121470255695856
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.0006679432326305856

$ws.Range("B7").Value = 0.3434965348879765
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.002132157156258965
$ws.Range("E7").Value = -0.00039249942726090606
$ws.Range("F7").Value = -0.0020097794672152285
$ws.Range("G7").Value = 0.0010487169578024236
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.0001871114577394685
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.002425503416774666

$ws.Range("B8").Value = 0.1788441461247557
$ws.Range("C8").Value = -0.04265844732882776
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.00006414713462962209
$ws.Range("F8").Value = -0.0004988726455031549
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.00004305772349578166
$ws.Range("I8").Value = 0.001610932322917725
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = -0.00005495409465328116

$ws.Range("B9").Value = 0.20918921556163528
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.0015911595810978302
$ws.Range("E9").Value = -0.002826774000271437
$ws.Range("F9").Value = -0.006863063338547238
$ws.Range("G9").Value = 0.000933775063616694
$ws.Range("H9").Value = -0.00017961081608560976
$ws.Range("I9").Value = -0.0004335672681842187
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.000355315079281493

$ws.Range("B10").Value = 0.429788748700226
$ws.Range("C10").Value = 0.07540769026031996
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.000046545501517275485
$ws.Range("F10").Value = -0.0004305793641990907
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.0000033079128631505172
$ws.Range("I10").Value = -0.0005517884614217137
$ws.Range("J10").Value = -0.0022001050535813906
$ws.Range("K10").Value = 0.0008264722260445945

$ws.Range("B11").Value = 0.3355115014665242
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = -0.011661268035431337
$ws.Range("E11").Value = 0.001597754961947659
$ws.Range("F11").Value = 0.0022075682020081732
$ws.Range("G11").Value = 0.0018212067562167022
$ws.Range("H11").Value = 0.00026856595862815674
$ws.Range("I11").Value = -0.00410059009696332
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.010170409674632397

# --- Append new row 12 for 2025-08-30 ---
# Write the date through a text formula then paste as value so it lands as literal
# text (shared string) instead of being auto-converted to a date serial number,
# keeping the default (unstyled) cell format used throughout the rest of the sheet.
$ws.Range("Z1").Formula = "=""2025-08-30"""
$ws.Range("Z1").Copy()
$ws.Range("A12").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B12").Value = 0.3307987144539396
$ws.Range("C12").Value = -0.039247296765689596
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.00012784888326927176
$ws.Range("F12").Value = 0.000010774763087169364
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.000004168558105142924
$ws.Range("I12").Value = -0.0001319112995502396
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.0031650590149974733
